$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 4 (shifts the existing Uruguay row down to row 5)
$ws.Rows("4:4").Insert()

# Populate the newly inserted row 4 with the new Brazil Serie A match
$ws.Range("A4").Value = "8nJEo620"
$ws.Range("B4").Value = "26/11/2024"
$ws.Range("C4").Value = "20:00"
$ws.Range("D4").Value = "BRAZIL - SERIE A BETANO"
$ws.Range("E4").Value = "Fortaleza"
$ws.Range("F4").Value = "Flamengo RJ"
$ws.Range("G4").Value = 3.1
$ws.Range("H4").Value = 3.2
$ws.Range("I4").Value = 2.4
$ws.Range("J4").Value = 3.75
$ws.Range("K4").Value = 2.05
$ws.Range("L4").Value = 3.1
$ws.Range("M4").Value = 1.07
$ws.Range("N4").Value = 8.5
$ws.Range("O4").Value = 1.36
$ws.Range("P4").Value = 3.2
$ws.Range("Q4").Value = 2.1
$ws.Range("R4").Value = 1.7
$ws.Range("S4").Value = 1.44
$ws.Range("T4").Value = 2.63
$ws.Range("U4").Value = 1.91
$ws.Range("V4").Value = 1.91
$ws.Range("W4").Value = 9
$ws.Range("X4").Value = 15
$ws.Range("Y4").Value = 11
$ws.Range("Z4").Value = 34
$ws.Range("AA4").Value = 26
$ws.Range("AB4").Value = 34
$ws.Range("AC4").Value = 8.5
$ws.Range("AD4").Value = 6
$ws.Range("AE4").Value = 15
$ws.Range("AF4").Value = 51
$ws.Range("AG4").Value = 7.5
$ws.Range("AH4").Value = 11
$ws.Range("AI4").Value = 9.5
$ws.Range("AJ4").Value = 23
$ws.Range("AK4").Value = 21
$ws.Range("AL4").Value = 29
$ws.Range("AM4").Value = 301
$ws.Range("AN4").Value = 5
$ws.Range("AO4").Value = 17
$ws.Range("AP4").Value = 29
$ws.Range("AQ4").Value = 51
$ws.Range("AR4").Value = 81
$ws.Range("AS4").Value = 201
$ws.Range("AT4").Value = 2.63
$ws.Range("AU4").Value = 8
$ws.Range("AV4").Value = 51
$ws.Range("AW4").Value = 4.33
$ws.Range("AX4").Value = 13
$ws.Range("AY4").Value = 23
$ws.Range("AZ4").Value = 41
$ws.Range("BA4").Value = 67
$ws.Range("BB4").Value = 201
$ws.Range("BC4").Value = 126
$ws.Range("BD4").Value = 126

# Update a handful of odds on the shifted row (now row 5) that changed too
$ws.Range("M5").Value = 1.06
$ws.Range("N5").Value = 10
$ws.Range("S5").Value = 1.36
$ws.Range("T5").Value = 3
$ws.Range("AB5").Value = 29
$ws.Range("AC5").Value = 10
$ws.Range("AG5").Value = 9.5
$ws.Range("AJ5").Value = 29
$ws.Range("AT5").Value = 3
$ws.Range("AU5").Value = 7.5
$ws.Range("AY5").Value = 23
